# Update cryptocurrency price/volume data (and fix a row ordering swap for PEPE/TheGraph)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'66.880.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -5.62%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.219.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -8.87%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'579.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -5.48%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'151.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -13.00%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  -0.18%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'3.212.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -8.91%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  -11.19%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = "'  -12.52%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'6.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -7.56%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.500"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -15.21%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'38.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -17.73%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'0.0000243"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -11.94%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'3.740.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -8.87%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'66.824.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.73%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'543.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -11.66%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'3.219.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -8.87%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  -5.91%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'7.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -15.51%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'  -15.21%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "'  -14.66%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'7.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -13.95%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'85.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -12.70%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'13.46"
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.08%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'3.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -16.90%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'8.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -11.90%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'29.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -13.39%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'2.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -19.05%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  -15.02%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'  -13.05%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'541.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -11.30%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'6.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -19.84%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "'  -16.78%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  +0.09%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'53.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -7.00%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.0431"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -9.55%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = "'  -15.93%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'9.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -15.68%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.124"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -14.22%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'2.932.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -13.10%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'2.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -27.20%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.261"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -16.78%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "'0.0₃0585"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -20.97%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'2.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -19.05%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'  -0.04%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'25.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -19.76%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  -18.15%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'  -13.33%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'123.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -7.85%  "
$ws.Range("E51").Style = "Normal"

